$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 503. This shifts the existing rows 503:539
# down to 504:540 (and carries the date-number-format style on column D
# along with them), matching the rest of the weekly observations already
# in the sheet.
$ws.Rows("503:503").Insert()

# Populate the newly-inserted row 503 with this week's observation. The
# descriptive columns (market/region/category/quality/origin/classification)
# are identical for every row in this sheet, so copy them from the row
# right below (the old row 503, now shifted to row 504).
$ws.Range("A503").Value = 9
$ws.Range("B503").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C503").Value = "Metropolitana"
$ws.Range("D503").Value = 45021
$ws.Range("E503").Value = 13
$ws.Range("F503").Value = 100112044
$ws.Range("G503").Value = "Perejil"
$ws.Range("H503").Value = "Sin especificar"
$ws.Range("I503").Value = "Primera"
$ws.Range("J503").Value = 95
$ws.Range("K503").Value = 13000
$ws.Range("L503").Value = 14000
$ws.Range("M503").Value = 13526
$ws.Range("N503").Value = "$/docena de atados"
$ws.Range("O503").Value = "Región Metropolitana"
$ws.Range("P503").Value = 4509
$ws.Range("Q503").Value = 3
$ws.Range("R503").Value = "Hortaliza"
